# Fixing Error with numbers in last name
# Adds the missing student (ID 573910, Thomas Carrie, Home School) as a new
# row 206 on both "RawData" and "ScheduleData", pushing the former last
# row (987654 / a / b) down to row 207, and corrects the rotation-slot
# (ROT) columns for a batch of previously-misassigned "2 Rotations:
# Global, SMCS" students on the ScheduleData sheet.

$wb = $excel.ActiveWorkbook

$wsRaw = $wb.Worksheets.Item("RawData")
$wsSched = $wb.Worksheets.Item("ScheduleData")

# ---------------------------------------------------------------------
# 1) RawData: insert the new student before the old last row (987654)
# ---------------------------------------------------------------------
$wsRaw.Rows.Item(206).Insert()

$wsRaw.Cells.Item(206, 1).Value = 573910
$wsRaw.Cells.Item(206, 2).Value = "Thomas"
$wsRaw.Cells.Item(206, 3).Value = "Carrie"
$wsRaw.Cells.Item(206, 4).Value = "Home School"
$wsRaw.Cells.Item(206, 5).Value = "Y"
$wsRaw.Cells.Item(206, 6).Value = "Y"
$wsRaw.Cells.Item(206, 7).Value = "Y"

# ---------------------------------------------------------------------
# 2) ScheduleData: insert the matching computed row in the same spot
# ---------------------------------------------------------------------
$wsSched.Rows.Item(206).Insert()

$wsSched.Cells.Item(206, 1).Value = "3 Rotations: Global, Humanities, SMCS"
$wsSched.Cells.Item(206, 2).Value = 573910
$wsSched.Cells.Item(206, 3).Value = "Thomas"
$wsSched.Cells.Item(206, 4).Value = "Carrie"
$wsSched.Cells.Item(206, 5).Value = "GL"
$wsSched.Cells.Item(206, 6).Value = "H"
$wsSched.Cells.Item(206, 7).Value = "S"
$wsSched.Cells.Item(206, 8).Value = "GE"

# ---------------------------------------------------------------------
# 3) ScheduleData: fix the ROT1/ROT2/ROT3 slot assignment for the
#    "2 Rotations: Global, SMCS" rows that shifted around the new entry
# ---------------------------------------------------------------------
function Set-Row($row, $e, $f, $g) {
    $wsSched.Cells.Item($row, 5).Value = $e
    $wsSched.Cells.Item($row, 6).Value = $f
    $wsSched.Cells.Item($row, 7).Value = $g
}

Set-Row 173 "GE" "S"  "GL"
Set-Row 176 "S"  "GE" "GL"
Set-Row 177 "GL" "GE" "S"
Set-Row 179 "GE" "GL" "S"
Set-Row 181 "S"  "GE" "GL"
Set-Row 183 "GL" "GE" "S"
Set-Row 188 "GE" "S"  "GL"
Set-Row 189 "S"  "GE" "GL"
Set-Row 191 "GL" "GE" "S"
Set-Row 193 "GE" "GL" "S"
Set-Row 195 "S"  "GE" "GL"
Set-Row 196 "GL" "GE" "S"
Set-Row 198 "GE" "S"  "GL"
Set-Row 199 "S"  "GE" "GL"
Set-Row 200 "GL" "GE" "S"
